$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D2:E51 so numeric-looking strings (e.g. "113.85")
# are not auto-converted to numbers by Excel when we set .Value below.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "45.513.71"
$ws.Range("E2").Value = "  +6.60%  "
$ws.Range("D3").Value = "2.390.56"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "113.85"
$ws.Range("E5").Value = "  +8.59%  "
$ws.Range("D6").Value = "319.30"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  +3.43%  "
$ws.Range("D10").Value = "42.31"
$ws.Range("E10").Value = "  +6.72%  "
$ws.Range("D11").Value = "0.0931"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  +5.00%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "15.87"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "2.749.77"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "2.389.85"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").Value = "45.428.76"
$ws.Range("E18").Value = "  +6.06%  "
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "74.83"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").Value = "3.57"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").Value = "264.16"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("E25").Value = "  +6.12%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "7.79"
$ws.Range("E27").Value = "  +5.40%  "
$ws.Range("D28").Value = "11.31"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").Value = "2.36"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "39.51"
$ws.Range("E30").Value = "  +8.59%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "22.78"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0979"
$ws.Range("E32").Value = "  +13.72%  "
$ws.Range("D33").Value = "172.39"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  +10.33%  "
$ws.Range("D35").Value = "0.132"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "4.95"
$ws.Range("E36").Value = "  +8.79%  "
$ws.Range("E37").Value = "  +6.64%  "
$ws.Range("D38").Value = "4.16"
$ws.Range("E38").Value = "  +14.89%  "
$ws.Range("E39").Value = "  +9.26%  "
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  +11.91%  "
$ws.Range("E42").Value = "  +6.54%  "
$ws.Range("D43").Value = "13.70"
$ws.Range("E43").Value = "  +12.35%  "
$ws.Range("D44").Value = "100.04"
$ws.Range("E44").Value = "  -9.18%  "
$ws.Range("D45").Value = "71.77"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").Value = "88.34"
$ws.Range("E46").Value = "  +14.34%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "5.81"
$ws.Range("E48").Value = "  +12.89%  "
$ws.Range("D49").Value = "116.10"
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("E50").Value = "  +9.61%  "
$ws.Range("D51").Value = "1.59"
$ws.Range("E51").Value = "  +10.86%  "

# Restore default (Normal) style on the range so no stray text-format
# style index lingers on these cells (matches original formatting).
$ws.Range("D2:E51").Style = "Normal"
